$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.303030303030303
$ws.Range("C2").Value = 0.330508474576271
$ws.Range("D2").Value = 0.275641025641026
$ws.Range("E2").Value = 0.256637168141593
$ws.Range("F2").Value = 0.274223034734918

$ws.Range("B3").Value = 0.502164502164502
$ws.Range("C3").Value = 0.572033898305085
$ws.Range("D3").Value = 0.57051282051282
$ws.Range("E3").Value = 0.539823008849557
$ws.Range("F3").Value = 0.404936014625229

$ws.Range("B4").Value = 0.155844155844156
$ws.Range("C4").Value = 0.152542372881356
$ws.Range("D4").Value = 0.198717948717949
$ws.Range("E4").Value = 0.123893805309735
$ws.Range("F4").Value = 0.158135283363803

$ws.Range("B5").Value = 0.303030303030303
$ws.Range("C5").Value = 0.322033898305085
$ws.Range("D5").Value = 0.198717948717949
$ws.Range("E5").Value = 0.348082595870207
$ws.Range("F5").Value = 0.448811700182815

$ws.Range("B6").Value = 0.515151515151515
$ws.Range("C6").Value = 0.601694915254237
$ws.Range("D6").Value = 0.673076923076923
$ws.Range("E6").Value = 0.690265486725664
$ws.Range("F6").Value = 0.648994515539305
